$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the first bullet of the "Student" use-case block (currently
# "Ik wil het aanbod van alle campi in antwerpen bekijken") and insert
# a brand new bullet paragraph in front of it containing the new,
# first use case "Ik wil nieuws over de stad lezen". The _GoBack
# bookmark (originally sitting at the end of the last bullet "Ik wil
# zoeken naar contactgegevens") must move to sit at the end of this
# new first bullet instead.
# ------------------------------------------------------------------

$findRange = $d.Content
$found = $findRange.Find.Execute(
    "Ik wil het aanbod van alle campi in antwerpen bekijken",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the anchor bullet paragraph."
}

# Insert an empty paragraph immediately before the found one; Word
# clones the paragraph/run formatting of the following paragraph.
$findRange.InsertParagraphBefore()

# A temporary sentinel character is appended after the new bullet's
# text so that the position where the bookmark must end up (the very
# end of the bullet's text, right before the paragraph mark) is not
# itself a paragraph-mark boundary while we create the bookmark. That
# boundary position is mishandled by Bookmarks.Add directly, so we
# add the bookmark next to a throw-away character and then delete
# that character, leaving the bookmark correctly collapsed at the
# end of the real text.
$newPara = $d.Paragraphs.Item(28)
$newPara.Range.Text = "Ik wil nieuws over de stad lezenX"

$newPara = $d.Paragraphs.Item(28)
$bmPos = $newPara.Range.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$newPara = $d.Paragraphs.Item(28)
$sentinelStart = $newPara.Range.End - 2
$sentinelRange = $d.Range($sentinelStart, $sentinelStart + 1)
$sentinelRange.Text = ""
